$d = $word.ActiveDocument
$d.Content.Find.Execute("MONTAGEM DE CARGA", $false, $false, $false, $false, $false, $true, 1, $false, "^t", 2)
Write-Output "Found: $($d.Content.Find.Found)"
